$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '43.071.63'
$ws.Range("E2").Value = '  +0.84%  '
Set-TextValue $ws.Range("D3") '2.304.76'
$ws.Range("E3").Value = '  +0.71%  '
$ws.Range("E4").Value = '  -0.01%  '
Set-TextValue $ws.Range("D5") '300.74'
$ws.Range("E5").Value = '  -0.01%  '
Set-TextValue $ws.Range("D6") '98.02'
$ws.Range("E6").Value = '  -1.12%  '
$ws.Range("E7").Value = '  +4.06%  '
$ws.Range("E8").Value = '  -0.02%  '
Set-TextValue $ws.Range("D9") '0.517'
$ws.Range("E9").Value = '  +1.27%  '
Set-TextValue $ws.Range("D10") '35.61'
$ws.Range("E10").Value = '  -0.60%  '
Set-TextValue $ws.Range("D11") '0.0789'
$ws.Range("E11").Value = '  -0.02%  '
$ws.Range("E12").Value = '  -0.18%  '
Set-TextValue $ws.Range("D13") '17.92'
$ws.Range("E13").Value = '  +0.39%  '
Set-TextValue $ws.Range("D14") '6.87'
$ws.Range("E14").Value = '  +0.96%  '
Set-TextValue $ws.Range("D15") '2.663.29'
$ws.Range("E15").Value = '  +0.59%  '
Set-TextValue $ws.Range("D16") '2.316.48'
$ws.Range("E16").Value = '  +0.52%  '
Set-TextValue $ws.Range("D17") '0.785'
$ws.Range("E17").Value = '  -1.65%  '
Set-TextValue $ws.Range("D18") '42.970.84'
$ws.Range("E18").Value = '  +0.78%  '
Set-TextValue $ws.Range("D19") '13.41'
$ws.Range("E19").Value = '  +7.78%  '
Set-TextValue $ws.Range("D20") '0.0₃0906'
$ws.Range("E20").Value = '  +1.03%  '
Set-TextValue $ws.Range("D21") '6.11'
$ws.Range("E21").Value = '  -1.08%  '
Set-TextValue $ws.Range("D22") '68.19'
$ws.Range("E22").Value = '  +0.64%  '
Set-TextValue $ws.Range("D23") '239.02'
$ws.Range("E23").Value = '  +1.68%  '
Set-TextValue $ws.Range("D24") '2.19'
$ws.Range("E24").Value = '  -0.60%  '
Set-TextValue $ws.Range("D25") '0.998'
$ws.Range("E25").Value = '  -0.27%  '
Set-TextValue $ws.Range("D26") '2.42'
$ws.Range("E26").Value = '  -0.68%  '
Set-TextValue $ws.Range("D27") '24.70'
$ws.Range("E27").Value = '  +0.65%  '
Set-TextValue $ws.Range("D28") '168.26'
$ws.Range("E28").Value = '  +0.17%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range("D29") '2.04'
$ws.Range("E29").Value = '  -6.97%  '
$ws.Range("B30").Value = 'Cosmos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range("D30") '9.14'
$ws.Range("E30").Value = '  -0.12%  '
Set-TextValue $ws.Range("D31") '33.05'
$ws.Range("E31").Value = '  -3.82%  '
Set-TextValue $ws.Range("D32") '5.19'
$ws.Range("E32").Value = '  +4.43%  '
Set-TextValue $ws.Range("D33") '0.999'
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("B34").Value = 'Celestia'
$ws.Range("C34").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue $ws.Range("D34") '18.15'
$ws.Range("E34").Value = '  +4.15%  '
$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D35") '4.80'
$ws.Range("E35").Value = '  +4.41%  '
$ws.Range("E36").Value = '  -0.24%  '
Set-TextValue $ws.Range("D37") '0.0689'
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("E38").Value = '  +1.16%  '
$ws.Range("E39").Value = '  +1.07%  '
Set-TextValue $ws.Range("D40") '0.112'
$ws.Range("E41").Value = '  -2.59%  '
Set-TextValue $ws.Range("D42") '2.005.49'
$ws.Range("E42").Value = '  +0.89%  '
$ws.Range("E43").Value = '  +0.11%  '
Set-TextValue $ws.Range("D44") '2.15'
$ws.Range("E44").Value = '  -3.28%  '
Set-TextValue $ws.Range("D45") '10.15'
$ws.Range("E45").Value = '  +0.72%  '
Set-TextValue $ws.Range("D46") '17.49'
$ws.Range("E46").Value = '  -0.16%  '
$ws.Range("E47").Value = '  -1.96%  '
Set-TextValue $ws.Range("D48") '54.40'
$ws.Range("E48").Value = '  -2.02%  '
Set-TextValue $ws.Range("D49") '2.531.00'
$ws.Range("E49").Value = '  +0.75%  '
$ws.Range("E50").Value = '  +0.78%  '
Set-TextValue $ws.Range("D51") '73.44'
$ws.Range("E51").Value = '  +5.15%  '
